$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "https://codeforces.com/contest/1285/problem/C"
$ws.Range("C8").Value = "LCM property use hoga"

$ws.Range("C8").Select()
